# Fixing waa (weight-at-age) indexing for 2 sex
# - Correct the Sex=2 (row 3) weight-at-age values on Weight_At_Age (previously
#   a duplicate of the Sex=1 row) to the distinct, correct second-sex series.
# - Apply a black-font style to a block of cells (D8:D37) on Weight_At_Age
#   (selected while reviewing, no values entered) which mints a new font /
#   cell style in the workbook.
# - Leave the workbook with Maturity_At_Age as the active sheet/tab, with
#   cell E30 selected there (the last sheet the author was looking at), and
#   Weight_At_Age's own selection parked at J17.

$wb = $excel.ActiveWorkbook

$wsWaa = $wb.Worksheets.Item("Weight_At_Age")
$wsWaa.Activate()

# Corrected Sex = 2 weight-at-age values (row 3, columns D:AG)
$newValues = [ordered]@{
    "D3" = 1.1085
    "E3" = 1.4285000000000001
    "F3" = 1.7228000000000001
    "G3" = 1.9837
    "H3" = 2.2088999999999999
    "I3" = 2.3995000000000002
    "J3" = 2.5586000000000002
    "K3" = 2.6899000000000002
    "L3" = 2.7974000000000001
    "M3" = 2.8847999999999998
    "N3" = 2.9554999999999998
    "O3" = 3.0125000000000002
    "P3" = 3.0583999999999998
    "Q3" = 3.0951
    "R3" = 3.1244999999999998
    "S3" = 3.1480000000000001
    "T3" = 3.1667999999999998
    "U3" = 3.1817000000000002
    "V3" = 3.1936
    "W3" = 3.2031000000000001
    "X3" = 3.2107000000000001
    "Y3" = 3.2166999999999999
    "Z3" = 3.2214999999999998
    "AA3" = 3.2252999999999998
    "AB3" = 3.2282999999999999
    "AC3" = 3.2307000000000001
    "AD3" = 3.2326000000000001
    "AE3" = 3.2341000000000002
    "AF3" = 3.2353000000000001
    "AG3" = 3.2381000000000002
}

foreach ($addr in $newValues.Keys) {
    $wsWaa.Range($addr).Value = $newValues[$addr]
}

# Set a black font on column D, rows 8-37 (stray formatting left behind while
# the author scrolled/selected the new indexing block - no cell values).
for ($r = 8; $r -le 37; $r++) {
    $wsWaa.Range("D$r").Font.Color = 0
}

# Author's selection on Weight_At_Age ends up at J17
$wsWaa.Range("J17").Select()

# Author finishes on Maturity_At_Age, with E30 selected there
$wsMaturity = $wb.Worksheets.Item("Maturity_At_Age")
$wsMaturity.Activate()
$wsMaturity.Range("E30").Select()

Write-Output "Done"
